# Updates the cryptos list (Price / Volume(1h) columns) with refreshed
# values, matching the upstream GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.207.62"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "1.703.32"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.56"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5294"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2640"
$ws.Range("E8").Value = "  -4.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06559"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.71"
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07635"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.567"
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").Value = "1.724.84"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "1.939.42"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5714"
$ws.Range("E15").Value = "  -4.81%  "
$ws.Range("D16").Value = "0.0₅8157"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.33"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "27.186.01"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.85"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.650"
$ws.Range("E21").Value = "  -3.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.41"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.951"
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.05"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.746"
$ws.Range("E26").Value = "  +6.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1216"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.239"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.26"
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05352"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.287"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("E32").Value = "  -6.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.403"
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.631"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.872"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9442"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5841"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01627"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.855"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "1.037.62"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.82"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").Value = "1.846.86"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.84"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4491"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06586"
$ws.Range("E50").Value = "  +10.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.046"
$ws.Range("E51").Value = "  -2.90%  "
